# DEQM Capability Statement Producer Client - QA pass:
#   - change actor names in text and figures, cap statements
#
# Concretely:
#   1. "profiles" sheet: remove the "Measure" and "Library" profile rows
#      (these are captured as separate CapabilityStatements elsewhere and
#      no longer belong in the Producer Client profile list).
#   2. "igs" sheet: point the "QI Core" row at its ImplementationGuide
#      canonical (instead of the old index.html doc page), and add a new
#      "CFQM" implementation guide row, styled like the other canonical
#      URIs (small monospace, dark blue).
#   3. Leave the "profiles" sheet active/selected (row 3 selected, as if
#      about to review/delete it) and update the "ops" sheet's lingering
#      selection.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. profiles sheet - drop the Measure and Library profile rows
# ---------------------------------------------------------------------
$profiles = $wb.Worksheets.Item("profiles")

# Row 3 = .../StructureDefinition/measure-deqm ("DEQM Measure Profile")
$profiles.Rows.Item(3).Delete()
# After the row-3 delete, the Library profile row (was row 6) is now row 5
$profiles.Rows.Item(5).Delete()

# ---------------------------------------------------------------------
# 2. igs sheet - update QI Core uri, add CFQM row
# ---------------------------------------------------------------------
$igs = $wb.Worksheets.Item("igs")

$igs.Range("B2").Value = "http://hl7.org/fhir/us/qicore/ImplementationGuide/qicore"

$igs.Range("A3").Value = "CFQM"
$igs.Range("B3").Value = "http://hl7.org/fhir/us/cqfmeasures/ImplementationGuide/cqfmeasures"
$igs.Range("B3").Font.Name = "Consolas"
$igs.Range("B3").Font.Size = 9
$igs.Range("B3").Font.Color = 6434563

# widen column B so the long canonical urls fit (was 32, now fits the
# longest uri at ~69 chars)
$igs.Columns.Item(2).ColumnWidth = 68.17

# ---------------------------------------------------------------------
# 3. view-state bookkeeping - match what Excel leaves behind once the
#    profiles sheet is the one being reviewed/edited
# ---------------------------------------------------------------------
$ops = $wb.Worksheets.Item("ops")
$ops.Range("B16").Select()

$igs.Range("D7").Select()

$profiles.Activate()
$profiles.Rows.Item(3).Select()
